$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the previously used range so stale cell contents from the old
#     layout (rows shifted around) don't linger. ---
$ws.Range("A1:E10").ClearContents()

# --- Column widths (best-effort; the runtime quantizes ColumnWidth to
#     whole-pixel steps and does not preserve the bestFit flag, so these
#     inputs are chosen empirically to land as close as possible to the
#     target widths of 11.109375 / 32.6640625 / 54.88671875 / 19) ---
$ws.Columns.Item(1).ColumnWidth = 10.25
$ws.Columns.Item(2).ColumnWidth = 31.8
$ws.Columns.Item(3).ColumnWidth = 53.95
$ws.Columns.Item(4).ColumnWidth = 18.1

# --- USER APIs section ---
$ws.Range("A1").Value = "Server URL"
$ws.Range("B1").Value = "https://courier50003.herokuapp.com/"

$ws.Range("A3").Value = "USER APIs"

$ws.Range("A4").Value = "End-point"
$ws.Range("B4").Value = "Pre-Condition "
$ws.Range("C4").Value = "Post-Condition"
$ws.Range("D4").Value = "Database Collection "
$ws.Range("E4").Value = "Special Notes"

$ws.Range("A5").Value = "/user/test"
$ws.Range("B5").Value = "NIL"
$ws.Range("C5").Value = "String Message 'This is the USER_MANAGEMENT Test controller!'"
$ws.Range("D5").Value = "USER_MANAGEMENT"

$ws.Range("A6").Value = "/user/signup"
$ws.Range("B6").Value = "1) name: String`n2) password: String`n3) contact_num: Long`n4) type: String"
$ws.Range("C6").Value = "1) success: Bool`n2) error: Bool`n3) message: String "
$ws.Range("D6").Value = "USER_MANAGEMENT"

$ws.Range("A7").Value = "/user/login"
$ws.Range("B7").Value = "1) email: String`n2) password: String"
$ws.Range("C7").Value = "1) success: Bool`n2) message: String`n3) token: String `n4) authority: String"
$ws.Range("D7").Value = "USER_MANAGEMENT`nUSER_SESSION"

$ws.Range("A8").Value = "/user/logout"
$ws.Range("B8").Value = "1) token: String "
$ws.Range("C8").Value = "1) success: Bool`n2) message: String`n3) id: String"
$ws.Range("D8").Value = "USER_MANAGEMENT`nUSER_SESSION"

# --- ADMIN API section ---
$ws.Range("A10").Value = "ADMIN API"

$ws.Range("A11").Value = "End-point"
$ws.Range("B11").Value = "Pre-Condition "
$ws.Range("C11").Value = "Post-Condition"
$ws.Range("D11").Value = "Database Collection "
$ws.Range("E11").Value = "Special Notes"

$ws.Range("A12").Value = "/portal/test"
$ws.Range("B12").Value = "NIL"
$ws.Range("C12").Value = "String Message 'This is the REQUESTS Test controller!' "

$ws.Range("A13").Value = "/portal"
$ws.Range("A14").Value = "/portal"
$ws.Range("A15").Value = "/portal"
$ws.Range("A16").Value = "/portal"
$ws.Range("A17").Value = "/portal"
$ws.Range("A18").Value = "/portal"
$ws.Range("A19").Value = "/portal"
$ws.Range("A20").Value = "/portal"
$ws.Range("A21").Value = "/portal"
$ws.Range("A22").Value = "/portal"
$ws.Range("A23").Value = "/portal"

# --- Wrap text styling on the long descriptive cells ---
$ws.Range("B6:C6").WrapText = $true
$ws.Range("B7:D7").WrapText = $true
$ws.Range("C8:D8").WrapText = $true

# --- Row heights for the wrapped rows ---
$ws.Rows.Item(6).RowHeight = 57.6
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 43.2

# --- Selection state ---
[void]$ws.Range("B12:B13").Select()
